# Daily attendance processing - 2025-12-23 11:52:38
# Normalizes the "Recorded By" (column G) values on the active sheet so that
# the list of recorder identities is written in a consistent order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact "Recorded By" strings that need to be reordered, mapped to their
# normalized replacement.
$replacements = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
